$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in "Aashish Sort 5" timing data (rows 127-131, columns P:U)
# Row 127 - Trial 1
$ws.Range("P127").Value = 2
$ws.Range("Q127").Value = 11
$ws.Range("R127").Value = 164
$ws.Range("S127").Value = 7250
$ws.Range("T127").Value = "Unmeasureable"
$ws.Range("U127").Value = "Unmeasureable"

# Row 128 - Trial 2
$ws.Range("P128").Value = 1
$ws.Range("Q128").Value = 2
$ws.Range("R128").Value = 108
$ws.Range("S128").Value = 6997
$ws.Range("T128").Value = "Unmeasureable"
$ws.Range("U128").Value = "Unmeasureable"

# Row 129 - Trial 3
$ws.Range("P129").Value = 1
$ws.Range("Q129").Value = 3
$ws.Range("R129").Value = 90
$ws.Range("S129").Value = 7124
$ws.Range("T129").Value = "Unmeasureable"
$ws.Range("U129").Value = "Unmeasureable"

# Row 130 - Trial 4
$ws.Range("P130").Value = 1
$ws.Range("Q130").Value = 1
$ws.Range("R130").Value = 78
$ws.Range("S130").Value = 6951
$ws.Range("T130").Value = "Unmeasureable"
$ws.Range("U130").Value = "Unmeasureable"

# Row 131 - Trial 5
$ws.Range("P131").Value = 1
$ws.Range("Q131").Value = 1
$ws.Range("R131").Value = 72
$ws.Range("S131").Value = 6908
$ws.Range("T131").Value = "Unmeasureable"
$ws.Range("U131").Value = "Unmeasureable"

# Update the selection to match where the author left off
$ws.Range("V129").Select()
